# Fruta / hortaliza, semanal
# Insert 3 new weekly price records (dated 2023-03-10) into the
# "Macroferia Regional de Talca - Uva" consolidated sheet, right before the
# existing row that was previously row 553 (now pushed down to row 556).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 553-555 (everything from old row 553 onward
# shifts down by 3, old row 642 becomes new row 645).
$ws.Range("A553:A555").EntireRow.Insert()

function Set-Row($r, $vals) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

Set-Row 553 @(5,"Macroferia Regional de Talca","Maule",44995,7,"Fruta",100109,"Uva",100109001,"Uva","Flame Seedless","Primera",230,9000,9000,9000,"`$/bandeja 18 kilos","Región de O'Higgins",500,18)
Set-Row 554 @(5,"Macroferia Regional de Talca","Maule",44995,7,"Fruta",100109,"Uva",100109001,"Uva","Red Globe","Primera",260,9000,10000,9577,"`$/bandeja 18 kilos","Región de O'Higgins",532,18)
Set-Row 555 @(5,"Macroferia Regional de Talca","Maule",44995,7,"Fruta",100109,"Uva",100109001,"Uva","Superior Seedless","Primera",240,9000,9000,9000,"`$/bandeja 18 kilos","Región de O'Higgins",500,18)

# Match the date formatting used by the rest of the "Fecha" column.
$ws.Range("D553:D555").NumberFormat = "YYYY-MM-DD HH:MM:SS"
